# Create new excel template to match DB schema changes
# - clears the sample shipment rows (rows 2-5) leaving the header row intact
# - adds a left-hand thin border on column D (row 2 only getting the larger
#   "size 12" font treatment) and column J (rows 2-5)
# - moves the active selection to G9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump D2's font to match the size-12 font used elsewhere in the template
# before the border is applied, so the resulting style carries both the
# font and the border change.
$ws.Range("D2").Font.Size = 12

# Add a thin left border on column D (rows 2-5) and column J (rows 2-5).
$ws.Range("D2:D5").Borders.Item(7).LineStyle = 1
$ws.Range("J2:J5").Borders.Item(7).LineStyle = 1

# Wipe out the sample data values (rows 2-5, columns A-Q) while keeping
# the formatting/styles that were just applied (and any pre-existing ones
# such as the text-number formatting on columns I/Q).
$ws.Range("A2:Q5").ClearContents()

# Move the selection like the author left it.
$ws.Range("G9").Select()
